# إضافة حدث جديد في Card12 by HOSSAM at 2025-12-08 11:35:59
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card12")

# Row 15: the previously-empty cells B15:K15 get the literal text "nan",
# matching the rest of the table's convention for "no value".
$ws.Range("B15:K15").Value = "nan"

# Row 16: new service event row appended below row 15.
# A leading "'" forces text storage (so "12"/"11/5/2025" aren't
# auto-converted to a number/date, and blank cells stay text-typed empty
# strings like the rest of the sheet's inlineStr cells); re-applying the
# "Normal" style afterwards clears the quote-prefix marker Excel leaves
# behind so the cell format matches the untouched rows above.
$cols = 1..15
foreach ($col in $cols) {
    $cell = $ws.Cells.Item(16, $col)
    switch ($col) {
        1  { $cell.Value = "'12" }
        12 { $cell.Value = "'11/5/2025" }
        13 { $cell.Value = "فني" }
        14 { $cell.Value = "قطع سير كويلر مسنن 1270" }
        15 { $cell.Value = "تم تغير سير 1270" }
        default { $cell.Value = "'" }
    }
    if ($col -le 12) {
        $cell.Style = "Normal"
    }
}
